$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.144.82'
$ws.Range('E2').Value = '  -0.88%  '
$ws.Range('D3').Value = '2.263.93'
$ws.Range('E3').Value = '  -1.10%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.34'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.24%  '
$ws.Range('E7').Value = '  -0.97%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  -1.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.89'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0789'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.19%  '
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.90'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.24%  '
$ws.Range('D14').Value = '2.615.90'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.65'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.85%  '
$ws.Range('D16').Value = '2.259.83'
$ws.Range('E16').Value = '  -1.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.790'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.63%  '
$ws.Range('D18').Value = '42.019.79'
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.24'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.81%  '
$ws.Range('E20').Value = '  -2.16%  '
$ws.Range('E21').Value = '  -0.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.82'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.73%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.57'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.72%  '
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.46'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.54'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.18%  '
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '162.48'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.22'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.66%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.17'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.72%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.60'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.94%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0735'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.65%  '
$ws.Range('E37').Value = '  -0.96%  '
$ws.Range('E38').Value = '  -4.58%  '
$ws.Range('E39').Value = '  -1.03%  '
$ws.Range('E40').Value = '  -1.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.05'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.72%  '
$ws.Range('E42').Value = '  +2.75%  '
$ws.Range('D43').Value = '1.947.74'
$ws.Range('E43').Value = '  -3.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '18.94'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.09%  '
$ws.Range('E45').Value = '  -1.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.91'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.18%  '
$ws.Range('E47').Value = '  -3.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '53.88'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.78%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '91.93'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '71.47'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.70%  '
$ws.Range('E51').Value = '  -2.77%  '
